# Refresh the cryptocurrency symbol snapshot (Price + Volume(1h) columns)
# with the latest scrape. Columns D and E hold plain-text values (prices
# like "0.9340" and percentages like "-2.81%") rather than real numbers,
# so each cell is forced to Text format before the new value is written;
# this stops Excel from reinterpreting a numeric-looking/percent-looking
# string as a true number (which would silently reformat it, e.g.
# "0.9340" -> 0.934). The style is reset back to "Normal" immediately
# afterwards so only the cell's text content changes -- no new formatting
# is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "321.49" },
    @{ Cell = "E2"; Value = "-2.78%" },
    @{ Cell = "D3"; Value = "42.60" },
    @{ Cell = "E3"; Value = "-6.28%" },
    @{ Cell = "D4"; Value = "5.182" },
    @{ Cell = "E4"; Value = "-7.19%" },
    @{ Cell = "D5"; Value = "0.08183" },
    @{ Cell = "E5"; Value = "-1.92%" },
    @{ Cell = "D6"; Value = "4.312" },
    @{ Cell = "E6"; Value = "-3.08%" },
    @{ Cell = "D7"; Value = "1.814" },
    @{ Cell = "E7"; Value = "-13.47%" },
    @{ Cell = "D8"; Value = "0.9341" },
    @{ Cell = "E8"; Value = "-4.68%" },
    @{ Cell = "D9"; Value = "0.1110" },
    @{ Cell = "E9"; Value = "-7.83%" },
    @{ Cell = "D10"; Value = "0.1865" },
    @{ Cell = "E10"; Value = "-3.23%" },
    @{ Cell = "D11"; Value = "0.09454" },
    @{ Cell = "E11"; Value = "-4.14%" },
    @{ Cell = "D12"; Value = "0.04681" },
    @{ Cell = "E12"; Value = "0.24%" },
    @{ Cell = "D13"; Value = "7.413" },
    @{ Cell = "E13"; Value = "-28.23%" },
    @{ Cell = "D14"; Value = "0.1057" },
    @{ Cell = "E14"; Value = "-0.02%" },
    @{ Cell = "D15"; Value = "0.001308" },
    @{ Cell = "E15"; Value = "2.25%" },
    @{ Cell = "D16"; Value = "0.005745" },
    @{ Cell = "E16"; Value = "-2.57%" },
    @{ Cell = "E17"; Value = "-0.70%" },
    @{ Cell = "D18"; Value = "2.533" },
    @{ Cell = "E18"; Value = "-0.40%" },
    @{ Cell = "D19"; Value = "0.3380" },
    @{ Cell = "E19"; Value = "0.90%" },
    @{ Cell = "E20"; Value = "-0.06%" },
    @{ Cell = "E21"; Value = "-8.45%" },
    @{ Cell = "D22"; Value = "0.04162" },
    @{ Cell = "E22"; Value = "-0.34%" },
    @{ Cell = "D23"; Value = "0.001249" },
    @{ Cell = "E23"; Value = "-3.68%" },
    @{ Cell = "D24"; Value = "0.004330" },
    @{ Cell = "E24"; Value = "-5.03%" },
    @{ Cell = "D25"; Value = "0.0001200" },
    @{ Cell = "E25"; Value = "-7.86%" },
    @{ Cell = "D26"; Value = "0.0002980" },
    @{ Cell = "E26"; Value = "-20.48%" },
    @{ Cell = "D38"; Value = "0.02730" },
    @{ Cell = "E38"; Value = "1.23%" },
    @{ Cell = "D39"; Value = "0.05554" },
    @{ Cell = "E39"; Value = "-3.61%" },
    @{ Cell = "D40"; Value = "0.008014" },
    @{ Cell = "E40"; Value = "1.32%" },
    @{ Cell = "D41"; Value = "0.1398" },
    @{ Cell = "E41"; Value = "-2.40%" },
    @{ Cell = "D42"; Value = "0.006550" },
    @{ Cell = "E42"; Value = "-12.90%" },
    @{ Cell = "D43"; Value = "0.002065" },
    @{ Cell = "E43"; Value = "-1.65%" },
    @{ Cell = "D44"; Value = "0.008280" },
    @{ Cell = "E44"; Value = "-2.43%" },
    @{ Cell = "D45"; Value = "0.3494" },
    @{ Cell = "E45"; Value = "3.78%" },
    @{ Cell = "D46"; Value = "0.00006912" },
    @{ Cell = "E46"; Value = "-2.54%" },
    @{ Cell = "D47"; Value = "0.00000000750" },
    @{ Cell = "E47"; Value = "-0.16%" },
    @{ Cell = "D48"; Value = "0.003511" },
    @{ Cell = "E48"; Value = "-0.53%" },
    @{ Cell = "D49"; Value = "0.003531" },
    @{ Cell = "E49"; Value = "-0.17%" },
    @{ Cell = "D50"; Value = "0.00002101" },
    @{ Cell = "E50"; Value = "-0.16%" },
    @{ Cell = "D51"; Value = "0.0002001" },
    @{ Cell = "E51"; Value = "-0.16%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
